$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.887.83"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.842.91"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.19%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "309.42"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4750"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.32%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3666"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.11%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07192"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.9258"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  +1.88%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07687"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "1.861.94"
$ws.Range("E13").Value = "  +1.86%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.300"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.389"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "88.61"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("E17").Value = "  -0.09%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008614"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "26.904.25"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.54"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.95%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.047"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.62"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.05%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.918"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "152.23"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "18.12"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.46%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.999"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.82%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "114.16"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "4.926"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.08851"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.300"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +5.47%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.7485"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +3.06%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.171"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +4.41%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.477"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.23%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.709"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.090"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.01949"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.60%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.05259"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +3.14%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.959"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.46%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.5184"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +3.08%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "6.961"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.83%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1509"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.39%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.199"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.29%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.49"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +5.62%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.4722"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.88%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "101.60"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.80%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.598"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "65.38"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("E50").Value = "  +0.67%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.8871"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +4.68%  "
